# Apply the edits described by the commit diff to Sheet1:
#   - G37:G46 gain a new value of 300 (column G extended down to row 46)
#   - F51:F60 gain a new value of 300 (column F extended down to row 60)
#   - Ten new data rows (161-170) are appended, continuing the A/E series
#     (A = 159..168, E = 300 for each)
#   - The sheet's selection/view is left on G37:G46 (activeCell G37),
#     matching the saved selection in the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G values for rows 37-46 ---
$ws.Range("G37:G46").Value = 300

# --- New column F values for rows 51-60 ---
$ws.Range("F51:F60").Value = 300

# --- Append new rows 161-170 (A: 159..168, E: 300) ---
for ($i = 0; $i -lt 10; $i++) {
    $row = 161 + $i
    $ws.Cells.Item($row, 1).Value = 159 + $i
    $ws.Cells.Item($row, 5).Value = 300
}

# --- Restore the saved selection/view state ---
$ws.Activate() | Out-Null
$ws.Range("G37:G46").Select() | Out-Null
